# Auto-generated: apply per-cell numeric updates from the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1480.0714
$ws.Range("I28").Value = 226.46153
$ws.Range("J28").Value = 17777
$ws.Range("K28").Value = 226.46153
$ws.Range("L28").Value = 17777
$ws.Range("M28").Value = 258.53847
$ws.Range("N28").Value = -18747

$ws.Range("H62").Value = 1993.9375
$ws.Range("I62").Value = 1954.5
$ws.Range("J62").Value = 2112.25
$ws.Range("K62").Value = 1954.5
$ws.Range("L62").Value = 2112.25
$ws.Range("M62").Value = -1330.5
$ws.Range("N62").Value = -3360.25

$ws.Range("H65").Value = 1993.9375
$ws.Range("I65").Value = 1954.5
$ws.Range("J65").Value = 2112.25
$ws.Range("K65").Value = 9772.5
$ws.Range("L65").Value = 10561.25
$ws.Range("M65").Value = -6652.5
$ws.Range("N65").Value = -16801.25

$ws.Range("H107").Value = 757.05884
$ws.Range("I107").Value = 554
$ws.Range("J107").Value = 899.2
$ws.Range("K107").Value = 554
$ws.Range("L107").Value = 899.2
$ws.Range("M107").Value = 1366
$ws.Range("N107").Value = -4739.2

$ws.Range("H113").Value = 3166.7693
$ws.Range("J113").Value = 3296
$ws.Range("L113").Value = 3296
$ws.Range("N113").Value = -9804

$ws.Range("H129").Value = 17852.17
$ws.Range("J129").Value = 20590.96
$ws.Range("L129").Value = 61772.88
$ws.Range("N129").Value = -71772.88

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1833.6786
$ws.Range("I61").Value = 1406.8096
$ws.Range("J61").Value = 3114.2856
$ws.Range("K61").Value = 1406.8096
$ws.Range("L61").Value = 3114.2856
$ws.Range("M61").Value = -1194.8096
$ws.Range("N61").Value = -3538.2856

$ws.Range("H74").Value = 1895.1666
$ws.Range("I74").Value = 1567.579
$ws.Range("K74").Value = 1567.579
$ws.Range("M74").Value = -693.579

$ws.Range("H77").Value = 1895.1666
$ws.Range("I77").Value = 1567.579
$ws.Range("K77").Value = 7837.895
$ws.Range("M77").Value = -3469.895

$ws.Range("H122").Value = 1575.08
$ws.Range("I122").Value = 1159.2667
$ws.Range("K122").Value = 3477.800099999999
$ws.Range("M122").Value = -1027.800099999999

$ws.Range("H136").Value = 1833.6786
$ws.Range("I136").Value = 1406.8096
$ws.Range("J136").Value = 3114.2856
$ws.Range("K136").Value = 4220.4288
$ws.Range("L136").Value = 9342.856800000001
$ws.Range("M136").Value = -1670.4288
$ws.Range("N136").Value = -14442.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 333.89474
$ws.Range("I80").Value = 426
$ws.Range("J80").Value = 301
$ws.Range("K80").Value = 426
$ws.Range("L80").Value = 301
$ws.Range("M80").Value = 572
$ws.Range("N80").Value = -2297

$ws.Range("H83").Value = 333.89474
$ws.Range("I83").Value = 426
$ws.Range("J83").Value = 301
$ws.Range("K83").Value = 2130
$ws.Range("L83").Value = 1505
$ws.Range("M83").Value = 2862
$ws.Range("N83").Value = -11489

$ws.Range("H86").Value = 1532.0834
$ws.Range("I86").Value = 1368.5625
$ws.Range("J86").Value = 1859.125
$ws.Range("K86").Value = 1368.5625
$ws.Range("L86").Value = 1859.125
$ws.Range("M86").Value = -245.5625
$ws.Range("N86").Value = -4105.125

$ws.Range("H89").Value = 1532.0834
$ws.Range("I89").Value = 1368.5625
$ws.Range("J89").Value = 1859.125
$ws.Range("K89").Value = 6842.8125
$ws.Range("L89").Value = 9295.625
$ws.Range("M89").Value = -1226.8125
$ws.Range("N89").Value = -20527.625

$ws.Range("H134").Value = 68661.734
$ws.Range("I134").Value = 112434.664
$ws.Range("J134").Value = 3002.3333
$ws.Range("K134").Value = 337303.992
$ws.Range("L134").Value = 9006.999899999999
$ws.Range("M134").Value = -334768.992
$ws.Range("N134").Value = -14076.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 54.166668
$ws.Range("I7").Value = 36.4
$ws.Range("J7").Value = 76.375
$ws.Range("K7").Value = 36.4
$ws.Range("L7").Value = 76.375
$ws.Range("M7").Value = 76.59999999999999
$ws.Range("N7").Value = -302.375

$ws.Range("H16").Value = 878.125
$ws.Range("I16").Value = 849.75
$ws.Range("J16").Value = 906.5
$ws.Range("K16").Value = 849.75
$ws.Range("L16").Value = 906.5
$ws.Range("M16").Value = -562.75
$ws.Range("N16").Value = -1480.5

$ws.Range("H107").Value = 549.1053000000001
$ws.Range("I107").Value = 506.29413
$ws.Range("J107").Value = 913
$ws.Range("K107").Value = 506.29413
$ws.Range("L107").Value = 913
$ws.Range("M107").Value = 1413.70587
$ws.Range("N107").Value = -4753

$ws.Range("H113").Value = 878.125
$ws.Range("I113").Value = 849.75
$ws.Range("J113").Value = 906.5
$ws.Range("K113").Value = 849.75
$ws.Range("L113").Value = 906.5
$ws.Range("M113").Value = 1320.25
$ws.Range("N113").Value = -5246.5

$ws.Range("H134").Value = 838.7727
$ws.Range("I134").Value = 794.6667
$ws.Range("K134").Value = 2384.0001
$ws.Range("M134").Value = 150.9998999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 754.0476
$ws.Range("I5").Value = 490.76923
$ws.Range("J5").Value = 1181.875
$ws.Range("K5").Value = 1472.30769
$ws.Range("L5").Value = 3545.625
$ws.Range("M5").Value = -1360.30769
$ws.Range("N5").Value = -3769.625

$ws.Range("H129").Value = 1049.7142
$ws.Range("I129").Value = 838.3333
$ws.Range("J129").Value = 1208.25
$ws.Range("K129").Value = 2514.9999
$ws.Range("L129").Value = 3624.75
$ws.Range("M129").Value = 2485.0001
$ws.Range("N129").Value = -13624.75

$ws.Range("H131").Value = 4301
$ws.Range("I131").Value = 9571.637000000001
$ws.Range("J131").Value = 890.58826
$ws.Range("K131").Value = 28714.911
$ws.Range("L131").Value = 2671.76478
$ws.Range("M131").Value = -23674.911
$ws.Range("N131").Value = -12751.76478

$ws.Range("H135").Value = 754.0476
$ws.Range("I135").Value = 490.76923
$ws.Range("J135").Value = 1181.875
$ws.Range("K135").Value = 4416.92307
$ws.Range("L135").Value = 10636.875
$ws.Range("M135").Value = -1881.92307
$ws.Range("N135").Value = -15706.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 737.64703
$ws.Range("I107").Value = 661.6667
$ws.Range("J107").Value = 920
$ws.Range("K107").Value = 661.6667
$ws.Range("L107").Value = 920
$ws.Range("M107").Value = 1258.3333
$ws.Range("N107").Value = -4760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1417.8572
$ws.Range("I100").Value = 1404.1666
$ws.Range("K100").Value = 1404.1666
$ws.Range("M100").Value = -863.1666

$ws.Range("H122").Value = 3036.4375
$ws.Range("I122").Value = 3044.8462
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 9134.5386
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6684.5386
$ws.Range("N122").Value = -13900

$ws.Range("H132").Value = 3543.3794
$ws.Range("I132").Value = 3797.5625
$ws.Range("K132").Value = 11392.6875
$ws.Range("M132").Value = -8862.6875

$ws.Range("H139").Value = 35130.715
$ws.Range("I139").Value = 36000
$ws.Range("J139").Value = 34985.832
$ws.Range("K139").Value = 36000
$ws.Range("L139").Value = 34985.832
$ws.Range("M139").Value = -30860
$ws.Range("N139").Value = -45265.832

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2098.9285
$ws.Range("I132").Value = 1458.8
$ws.Range("J132").Value = 3699.25
$ws.Range("K132").Value = 4376.4
$ws.Range("L132").Value = 11097.75
$ws.Range("M132").Value = -1846.4
$ws.Range("N132").Value = -16157.75

Write-Host "Updated 193 cells"
